$d = $word.ActiveDocument

# Revert the "Platforms for communication" paragraph back to its original
# wording: collapse " | Everything will be done in Teams 😀 except we use
# Trello for project planning" down to " | Eveything will be done in Teams :D"
$old = " | Everything will be done in Teams 😀 except we use Trello for project planning"
$new = " | Eveything will be done in Teams :D"

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
